{"js": "// Replace each three-digit-by-one-digit multiplication expression with its\n// updated version. Every \"old\" value below occurs exactly once in the\n// document body, so an exact, case-sensitive whole-text search-and-replace\n// is safe and unambiguous.\nconst replacements = [\n  [\"950\u00d72=1900\", \"590\u00d75=2950\"],\n  [\"387\u00d75=1935\", \"648\u00d77=4536\"],\n  [\"115\u00d74=460\", \"637\u00d76=3822\"],\n  [\"555\u00d78=4440\", \"311\u00d78=2488\"],\n  [\"887\u00d79=7983\", \"339\u00d79=3051\"],\n  [\"466\u00d79=4194\", \"693\u00d76=4158\"],\n  [\"669\u00d78=5352\", \"424\u00d77=2968\"],\n  [\"261\u00d79=2349\", \"260\u00d77=1820\"],\n  [\"356\u00d73=1068\", \"416\u00d74=1664\"],\n  [\"753\u00d72=1506\", \"370\u00d79=3330\"],\n  [\"236\u00d74=944\", \"138\u00d77=966\"],\n  [\"356\u00d74=1424\", \"239\u00d72=478\"],\n  [\"910\u00d72=1820\", \"150\u00d76=900\"],\n  [\"302\u00d79=2718\", \"566\u00d77=3962\"],\n  [\"318\u00d75=1590\", \"767\u00d76=4602\"],\n  [\"853\u00d78=6824\", \"276\u00d79=2484\"],\n  [\"580\u00d73=1740\", \"854\u00d73=2562\"],\n  [\"963\u00d74=3852\", \"733\u00d72=1466\"],\n  [\"731\u00d79=6579\", \"921\u00d75=4605\"],\n  [\"779\u00d76=4674\", \"516\u00d78=4128\"],\n  [\"416\u00d78=3328\", \"259\u00d77=1813\"],\n  [\"169\u00d73=507\", \"378\u00d74=1512\"],\n  [\"982\u00d72=1964\", \"638\u00d76=3828\"],\n  [\"554\u00d77=3878\", \"293\u00d73=879\"],\n  [\"402\u00d74=1608\", \"622\u00d76=3732\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression with its\n# updated version. Every \"old\" value occurs exactly once in the document, so\n# a plain Find/Replace-all (scoped to the whole story) is unambiguous for\n# each pair.\n$replacements = @(\n    @{ Old = \"950\u00d72=1900\"; New = \"590\u00d75=2950\" },\n    @{ Old = \"387\u00d75=1935\"; New = \"648\u00d77=4536\" },\n    @{ Old = \"115\u00d74=460\"; New = \"637\u00d76=3822\" },\n    @{ Old = \"555\u00d78=4440\"; New = \"311\u00d78=2488\" },\n    @{ Old = \"887\u00d79=7983\"; New = \"339\u00d79=3051\" },\n    @{ Old = \"466\u00d79=4194\"; New = \"693\u00d76=4158\" },\n    @{ Old = \"669\u00d78=5352\"; New = \"424\u00d77=2968\" },\n    @{ Old = \"261\u00d79=2349\"; New = \"260\u00d77=1820\" },\n    @{ Old = \"356\u00d73=1068\"; New = \"416\u00d74=1664\" },\n    @{ Old = \"753\u00d72=1506\"; New = \"370\u00d79=3330\" },\n    @{ Old = \"236\u00d74=944\"; New = \"138\u00d77=966\" },\n    @{ Old = \"356\u00d74=1424\"; New = \"239\u00d72=478\" },\n    @{ Old = \"910\u00d72=1820\"; New = \"150\u00d76=900\" },\n    @{ Old = \"302\u00d79=2718\"; New = \"566\u00d77=3962\" },\n    @{ Old = \"318\u00d75=1590\"; New = \"767\u00d76=4602\" },\n    @{ Old = \"853\u00d78=6824\"; New = \"276\u00d79=2484\" },\n    @{ Old = \"580\u00d73=1740\"; New = \"854\u00d73=2562\" },\n    @{ Old = \"963\u00d74=3852\"; New = \"733\u00d72=1466\" },\n    @{ Old = \"731\u00d79=6579\"; New = \"921\u00d75=4605\" },\n    @{ Old = \"779\u00d76=4674\"; New = \"516\u00d78=4128\" },\n    @{ Old = \"416\u00d78=3328\"; New = \"259\u00d77=1813\" },\n    @{ Old = \"169\u00d73=507\"; New = \"378\u00d74=1512\" },\n    @{ Old = \"982\u00d72=1964\"; New = \"638\u00d76=3828\" },\n    @{ Old = \"554\u00d77=3878\"; New = \"293\u00d73=879\" },\n    @{ Old = \"402\u00d74=1608\"; New = \"622\u00d76=3732\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n"}
